# Update the "data" worksheet:
#  - Column E (選點原因 reason codes / color hex) values are renumbered in the
#    shared-string table; the only real text change is the hex code that used
#    to read "#7FFF00" which is now "#00FF00".
#  - The header cells D1/E1 are relabeled to "選點原因" / "選點原因Color".
#  - The sheet view no longer freezes a scrolled-down topLeftCell and instead
#    just remembers the last selected cell, I80.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$lastRow = 141
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    if ($cell.Text -eq "#7FFF00") {
        $cell.Value = "#00FF00"
    }
}

$ws.Cells.Item(1, 4).Value = "選點原因"
$ws.Cells.Item(1, 5).Value = "選點原因Color"

$ws.Range("I80").Select()
